# Added user ability to approve or reject incumbent transfers, fixed a bug
# where incoming transfer request displayed requests already responded to.
#
# Data change: Row 11 corresponds to the requirement
# "* As a customer, I can accept a money transfer from another account."
# It is now marked as DONE in the Status column (C11), which causes its
# point value (2) to be credited in column D, and the running Total in D13
# to increase accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "accept a money transfer" requirement as DONE.
$ws.Range("C11").Value = "DONE"

# Recalculate so dependent formulas (D11, D13) pick up the new value.
$excel.Calculate()

# Move the active selection to C12, matching the cursor position left
# behind after completing the edit on C11.
$ws.Range("C12").Select()
